# Add a new "Czech" worksheet (cloned from "Belgium", same layout/styles)
# and populate it with the Czech market data, then update tab/selection
# state so Czech becomes the active sheet.

$wb = $excel.ActiveWorkbook
$belgium = $wb.Worksheets.Item("Belgium")

# Duplicate "Belgium" right after itself to get an identical template
# (same styles, merged cells, column layout) for the new "Czech" sheet.
$belgium.Copy($null, $belgium)
$czech = $wb.Worksheets.Item($belgium.Index + 1)
$czech.Name = "Czech"

# Fill in the Czech-specific values.
$czech.Range("B2").Value = "Czech Market"
$czech.Range("B4").Value = "NGC-3477/T1731"

# The Czech tab uses slightly narrower columns than Germany/Belgium.
$czech.Columns.Item(2).ColumnWidth = 34.21875
$czech.Columns.Item(3).ColumnWidth = 13.77734375
$czech.Columns.Item(4).ColumnWidth = 15.88671875

# Belgium is no longer the focused tab: it keeps a plain full-table
# selection instead of its previous single-cell selection.
[void]$belgium.Range("A1:D20").Select()

# Czech becomes the selected/active tab, with B13 highlighted.
[void]$czech.Range("B13").Select()
$czech.Activate()
